$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the donation amount in D3 (this will trigger recalculation of
# dependent formula cells G6, G8, G10)
$ws.Range("D3").Value = 50.5

# Force a full recalculation so the formula result caches are refreshed
$excel.CalculateFull()

# Update the selected cell/active cell to C3 (matches the saved selection
# state in the sheetView)
$ws.Range("C3").Select()
